$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: shrink height ---
$ws.Rows.Item(12).RowHeight = 34

# --- Row 13: reset height to default, add compass legend label in E13 ---
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Range("E13").Value = "COMPAS { SO,SW, WE, NW, NO, NE, EA, SE }"

# --- Row 14: compass-table header labels (H,J,L,N already have column style 2) ---
$ws.Range("H14").Value = "SO"
$ws.Range("I14").Value = "SW"
$ws.Range("J14").Value = "WE"
$ws.Range("K14").Value = "NW"
$ws.Range("L14").Value = "NO"
$ws.Range("M14").Value = "NE"
$ws.Range("N14").Value = "EA"
$ws.Range("O14").Value = "SE"

# --- Row 15: compass-table index values 0..7 ---
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 4
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = 6
$ws.Range("O15").Value = 7

# --- Compass grid (rows 17-21): condition text helper cells + formulas ---
$ws.Range("M17").Value = "NO"
$ws.Range("AG17").Value = "D0 - O0 == 0"
$ws.Range("AH17").Value = "D1 - O1 > 0"
$ws.Range("V17").Formula = '="if (" & AG17 & " && " & AH17 & ") { comp = " & M17 & ";}"'

$ws.Range("K18").Value = "NW"
$ws.Range("O18").Value = "NE"
$ws.Range("AE18").Value = "D0 - O0 < 0"
$ws.Range("AF18").Value = "D1 - O1 > 0"
$ws.Range("AI18").Value = "D0 - O0 > 0"
$ws.Range("AJ18").Value = "D1 - O1 > 0"
$ws.Range("V18").Formula = '="if (" & AE18 & " && " & AF18 & ") { comp = " & K18 & ";}"'
$ws.Range("Z18").Formula = '="if (" & AI18 & " && " & AJ18 & ") { comp = " & O18 & ";}"'

$ws.Range("K19").Value = "WE"
$ws.Range("O19").Value = "EA"
$ws.Range("AE19").Value = "D0 - O0 < 0"
$ws.Range("AF19").Value = "D1 - O1 == 0"
$ws.Range("AI19").Value = "D0 - O0 > 0"
$ws.Range("AJ19").Value = "D1 - O1 == 0"

$ws.Range("K20").Value = "SW"
$ws.Range("O20").Value = "SE"
$ws.Range("AE20").Value = "D0 - O0 < 0"
$ws.Range("AF20").Value = "D1 - O1 < 0"
$ws.Range("AI20").Value = "D0 - O0 > 0"
$ws.Range("AJ20").Value = "D1 - O1 < 0"

# Shared formulas: V19:V20 (si=0) then Z19:Z20 (si=1) -- order matters for si numbering.
$ws.Range("V19:V20").Formula = '="if (" & AE19 & " && " & AF19 & ") { comp = " & K19 & ";}"'
$ws.Range("Z19:Z20").Formula = '="if (" & AI19 & " && " & AJ19 & ") { comp = " & O19 & ";}"'

$ws.Range("M21").Value = "SO"
$ws.Range("AG21").Value = "D0 - O0 == 0"
$ws.Range("AH21").Value = "D1 - O1 < 0"
$ws.Range("V21").Formula = '="if (" & AG21 & " && " & AH21 & ") { comp = " & M21 & ";}"'

# --- Row 23: stray styled blank cell ---
$ws.Range("Z23").Style = $ws.Range("W5").Style

# --- Rows 24-28: "if ( compas == X ) { ...; ...;}" formulas (double-space quirk before trailing '&') ---
$ws.Range("V24").Formula = '="if ( compas == " & M17 & " ) { " & AG17  & ";' + "`n" + '" & AH17 & ";}"'
$ws.Range("V25").Formula = '="if ( compas == " & K18 & " ) { " & AE18  & ";' + "`n" + '" & AF18 & ";}"'
$ws.Range("Z25").Formula = '="if ( compas == " & O18 & " ) { " & AI18  & ";' + "`n" + '" & AJ18 & ";}"'
$ws.Range("V26").Formula = '="if ( compas == " & K19 & " ) { " & AE19  & ";' + "`n" + '" & AF19 & ";}"'
$ws.Range("Z26").Formula = '="if ( compas == " & O19 & " ) { " & AI19  & ";' + "`n" + '" & AJ19 & ";}"'
$ws.Range("V27").Formula = '="if ( compas == " & K20 & " ) { " & AE20  & ";' + "`n" + '" & AF20 & ";}"'
$ws.Range("Z27").Formula = '="if ( compas == " & O20 & " ) { " & AI20  & ";' + "`n" + '" & AJ20 & ";}"'
$ws.Range("V28").Formula = '="if ( compas == " & M21 & " ) { " & AG21  & ";' + "`n" + '" & AH21 & ";}"'

# --- Rows 30-37: final compiled "if ( compas == ... )" literal text blocks ---
$ws.Range("V30").Value = "if ( compas == NO ) { D0 - O0 == 0;`nD1 - O1 > 0;}"
$ws.Range("V31").Value = "if ( compas == NW ) { D0 - O0 < 0;`nD1 - O1 > 0;}"
$ws.Range("V32").Value = "if ( compas == WE ) { D0 - O0 < 0;`nD1 - O1 == 0;}"
$ws.Range("V33").Value = "if ( compas == SW ) { D0 - O0 < 0;`nD1 - O1 < 0;}"
$ws.Range("V34").Value = "if ( compas == SO ) { D0 - O0 == 0;`nD1 - O1 < 0;}"
$ws.Range("V35").Value = "if ( compas == NE ) { D0 - O0 > 0;`nD1 - O1 > 0;}"
$ws.Range("V36").Value = "if ( compas == EA ) { D0 - O0 > 0;`nD1 - O1 == 0;}"
$ws.Range("V37").Value = "if ( compas == SE ) { D0 - O0 > 0;`nD1 - O1 < 0;}"
